# Lesson 1.1 Course Introduction — "changed DSS to DRC in L1.1"
#
# Slide 18 (SlideID 280), shape 2 ("Late Policy" body placeholder) has a
# bullet that says the accommodation office is "Disability Services" and a
# sub-bullet that abbreviates it as "DSS". The office was renamed to the
# "Disability Resource Center" (DRC), so:
#   - "Disability Services, " -> "Disability Resource Center, "
#   - "DSS Accommodations ..." -> "DRC Accommodations ..."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Bullet: "If you have an accommodation from Disability Services, ..."
# Replace just the "Services, " word (keeping the surrounding text/runs
# intact) so it becomes "...Disability Resource Center, you must request...".
$para6 = $tr.Paragraphs(6)
$text6 = $para6.Text
$idx6 = $text6.IndexOf("Services, ")
if ($idx6 -ge 0) {
    $target = $tr.Characters($para6.Start + $idx6, "Services, ".Length)
    $target.Text = "Resource Center, "
}

# --- Sub-bullet: "DSS Accommodations are usually NOT available for Group Assignments"
# Re-fetch the paragraph since the text length of the deck shifted above.
# Replace the whole line in two steps (through an unrelated placeholder)
# so the run stays a single clean run instead of being split on the
# common "D" prefix shared by "DSS" and "DRC".
$para7 = $tr.Paragraphs(7)
$para7.Text = "PLACEHOLDER_FOR_DSS_TO_DRC_RENAME"
$para7 = $tr.Paragraphs(7)
$para7.Text = "DRC Accommodations are usually NOT available for Group Assignments"
